$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to text format so numeric-looking price
# strings (e.g. "1.00", "0.999") are preserved exactly as text, matching
# the original inlineStr cell type, instead of being coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '54.037.33'
$ws.Range("E2").Value = '  +0.77%  '
$ws.Range("D3").Value = '2.254.74'
$ws.Range("E3").Value = '  +2.56%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").Value = '492.73'
$ws.Range("E5").Value = '  +1.23%  '
$ws.Range("D6").Value = '127.47'
$ws.Range("E6").Value = '  +2.06%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +2.88%  '
$ws.Range("E10").Value = '  +2.41%  '
$ws.Range("E11").Value = '  +3.04%  '
$ws.Range("D12").Value = '4.66'
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("D13").Value = '2.667.75'
$ws.Range("E13").Value = '  +2.97%  '
$ws.Range("D14").Value = '21.74'
$ws.Range("E14").Value = '  +2.99%  '
$ws.Range("D15").Value = '53.993.20'
$ws.Range("E15").Value = '  +0.70%  '
$ws.Range("E16").Value = '  +0.40%  '
$ws.Range("D17").Value = '2.258.53'
$ws.Range("E17").Value = '  +1.40%  '
$ws.Range("D19").Value = '4.07'
$ws.Range("E19").Value = '  +2.77%  '
$ws.Range("E20").Value = '  +4.07%  '
$ws.Range("D21").Value = '298.98'
$ws.Range("E21").Value = '  +1.37%  '
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("E23").Value = '  -1.99%  '
$ws.Range("D24").Value = '62.03'
$ws.Range("E24").Value = '  -1.30%  '
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  +0.28%  '
$ws.Range("D26").Value = '0.369'
$ws.Range("E26").Value = '  +0.84%  '
$ws.Range("D27").Value = '2.360.34'
$ws.Range("E27").Value = '  +1.67%  '
$ws.Range("E28").Value = '  +1.95%  '
$ws.Range("E29").Value = '  +0.43%  '
$ws.Range("D30").Value = '166.78'
$ws.Range("E30").Value = '  +1.01%  '
$ws.Range("D31").Value = '1.60'
$ws.Range("E31").Value = '  +1.09%  '
$ws.Range("D32").Value = '5.84'
$ws.Range("E32").Value = '  +2.45%  '
$ws.Range("B33").Value = 'PEPE'
$ws.Range("C33").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D33").Value = '0.0₃0674'
$ws.Range("E33").Value = '  +1.69%  '
$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("D35").Value = '0.996'
$ws.Range("E35").Value = '  +0.28%  '
$ws.Range("E36").Value = '  -0.47%  '
$ws.Range("D37").Value = '17.60'
$ws.Range("E37").Value = '  +1.71%  '
$ws.Range("D38").Value = '0.887'
$ws.Range("E38").Value = '  +6.49%  '
$ws.Range("E39").Value = '  +2.61%  '
$ws.Range("D40").Value = '3.65'
$ws.Range("E40").Value = '  +3.13%  '
$ws.Range("D41").Value = '35.72'
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '1.38'
$ws.Range("E42").Value = '  +1.68%  '
$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D43").Value = '0.370'
$ws.Range("E43").Value = '  +0.73%  '
$ws.Range("E44").Value = '  +2.03%  '
$ws.Range("D45").Value = '124.83'
$ws.Range("E45").Value = '  -1.17%  '
$ws.Range("D46").Value = '4.71'
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("E47").Value = '  +0.51%  '
$ws.Range("D48").Value = '0.540'
$ws.Range("E48").Value = '  +1.01%  '
$ws.Range("E49").Value = '  +2.25%  '
$ws.Range("D50").Value = '234.61'
$ws.Range("E50").Value = '  +1.29%  '
$ws.Range("E51").Value = '  +0.92%  '

# Restore the default cell style on column D so no stray explicit
# number-format style lingers on cells that did not have one before.
$ws.Range("D2:D51").Style = "Normal"
